$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows/columns that are no longer needed first.
$ws.Rows("3:5").Delete()
$ws.Columns("F").Delete()

# Row 2 becomes the "Saldo Inicial" (opening balance) row instead of the
# old "Despesa / EQUIPAMENTOS" entry.
$ws.Range("A2").Value = "Saldo Inicial"
$ws.Range("B2").Value = "Inicial"
$ws.Range("C2").Value = 0

# D2 becomes a real date value (2025-01-01, serial 45658) instead of the
# free-text "26/01/2025" string, with a date/time number format applied.
$ws.Range("D2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").Value = 45658

$ws.Range("E2").Value = 0
